$d = $word.ActiveDocument

# Remove the embedded OLE object (the chapter-7 image) from its paragraph.
# The object is represented in the Word OM as an EMBED field; deleting the
# field removes the run that hosts it, leaving the paragraph empty.
$f = $d.Fields.Item(1)
$f.Delete()

# The document's trailing "_GoBack" bookmark (left over from the last edit
# position) moves from the end of the document to the now-empty paragraph
# that used to hold the picture.
$b = $d.Bookmarks.Item("_GoBack")
$b.Delete()

$p7 = $d.Paragraphs.Item(7)
$d.Bookmarks.Add("_GoBack", $p7.Range)
